$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4.68
$ws.Range("A3").Value = 157.46
$ws.Range("A4").Value = 11.82
$ws.Range("A5").Value = 131.97
$ws.Range("A6").Value = 45.52
